# Generate Report for Handoff
# The 228b4934-faed-44a0-b362-1b99a5cea0b9.md file has finished its handback
# cycle and is now ready to be handed off again. Update the status on the
# Overview sheet as well as on the per-language (zh-cn / de-de) detail
# sheets, and record the new "Latest Handoff Datetime" for each language.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: update status for the first file row ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"

# --- zh-cn sheet: update status + latest handoff datetime ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B2").Value = "Ready for handoff"
$zhcn.Range("D2").Value = "2016-02-22 14:02:31"

# --- de-de sheet: update status + latest handoff datetime ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B2").Value = "Ready for handoff"
$dede.Range("D2").Value = "2016-02-22 14:02:43"
